$d = $word.ActiveDocument

# Remove the "vnpt.SiteAddress" mail-merge placeholder run that follows
# "Địa chỉ: " in the representative's address line. The whole match is
# deleted (replaced with an empty string) so the run collapses away.
$d.Content.Find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
